$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title change: "Featureliste v.1.1" -> "Featureliste v.1.2" ---
$ws.Range("A1").Value = "Featureliste v.1.2"

# --- New columns J (Arbeitspaket) and K (Status) on header row 2 ---
# Copy formatting from existing header cell (I2) so J2/K2 match the header style
$ws.Range("I2").Copy()
$ws.Range("J2:K2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J2").Value = "Arbeitspaket"
$ws.Range("K2").Value = "Status"

# --- Row 3 gets an example Arbeitspaket entry ---
$ws.Range("J3").Value = "z.B. 1"

# --- New "Status" column K, value 0 for every data row (3-24) ---
for ($r = 3; $r -le 24; $r++) {
    $ws.Cells.Item($r, 11).Value = 0
}

# --- Fill in previously empty Aufwand / Risiko / Priorität (F, G, H) for rows 13-24 ---
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 5

$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 8

$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 10

$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 10

$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 2

$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 10

$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 8

$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 10

$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 10

$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 10

$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 5

$ws.Range("F24").Value = 4
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 4

# --- Row 15/17/18 Typ changes from "Funktionalität" to "Design" ---
$ws.Range("B15").Value = "Design"
$ws.Range("B17").Value = "Design"
$ws.Range("B18").Value = "Design"

# --- Row 15 "Ziel" (C15) loses its centered-alignment styling, back to default ---
$ws.Range("C15").Style = "Standard"

# --- Row 24 description text tweak (removed "nächsten") ---
$ws.Range("E24").Value = "Nachdem für den Spielzug gevotet wurde, werden die Upvotes für den Zug im Spiel dargestellt."

# --- New column K width ---
$ws.Columns.Item(11).ColumnWidth = 7.5

# --- Selection matches author's last cursor position ---
$ws.Range("D15").Select()
